$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 324. This shifts the existing
# rows 324:372 down to 326:374, which matches the dimension growing from
# A1:T372 to A1:T374.
$ws.Range("324:325").Insert()

# New weekly entries (2021-10-05, serial 44474) for "Pintón" and
# "Primera Pintón" qualities, inserted ahead of the existing history.
$newRows = @(
    @(5, "Macroferia Regional de Talca", "Maule", 44474, 7, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Pintón", 800, 18000, 18000, 18000, "`$/caja 20 kilos", "Ecuador", 900, 20),
    @(5, "Macroferia Regional de Talca", "Maule", 44474, 7, "Fruta", 100108, "Tropicales y subtropicales", 100108006, "Plátano", "Sin especificar", "Primera Pintón", 500, 19000, 19000, 19000, "`$/caja 20 kilos", "Ecuador", 950, 20)
)

$startRow = 324
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}

# Column D (Fecha) keeps the same datetime number format used throughout
# the rest of the column.
$ws.Range("D324:D325").NumberFormat = "YYYY-MM-DD HH:MM:SS"
